# issue #5: add legislator_id, name, date into dataframe
#
# The 股票 (stocks) sheet (4th worksheet) gets three new trailing columns:
#   H: date             -> "2012-04-16"
#   I: legislator_name  -> "呂學樟"
#   J: legislator_id    -> 892
#
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# --- Header row (row 1): copy the header formatting (bold + border, as used
# by the existing G1 header cell) onto the new header cells, then set text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("G1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("G1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# --- Data row (row 2): copy the plain data formatting from G2 onto the new
# data cells, then set values.
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$ws.Range("G2").Copy()
$ws.Range("J2").PasteSpecial(-4122)

# Force the date cell to be stored as literal text, not an auto-converted
# Excel date serial number.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "2012-04-16"

$ws.Range("I2").Value = "呂學樟"
$ws.Range("J2").Value = 892

$excel.CutCopyMode = $false
